# mydf.key.xlsx was re-created by opening mydf.key and keySave as xlsx.
# In the new export the "varlab" sheet/table no longer exists, and the
# "key" sheet's recode columns (value_old / value_new) gained a trailing
# ".", reflecting an extra "other/NA" level appended to each recode map.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the "varlab" worksheet (and, with it, its Table4) entirely.
$varlab = $wb.Worksheets.Item("varlab")
$varlab.Delete()

# Update the recode strings on the remaining "key" sheet to append "<."
# (ordered factors) / "|." (delimited lists) for the new trailing level.
$key = $wb.Worksheets.Item("key")

$key.Range("E4").Value = "lo<med<hi<."
$key.Range("F4").Value = "lo<mid<mid<."

$key.Range("E5").Value = "f|d|c|b|a|."
$key.Range("F5").Value = "f<d<c<b<a<."

$key.Range("E6").Value = "cindy|bobby|peter|marcia|greg|."
$key.Range("F6").Value = "Cindy<Bobby<Peter<Marcia<Greg<."

$key.Range("F7").Value = "fail<fail<pass<pass<pass<."

$key.Range("E8").Value = "1|2|3|4|5|."
$key.Range("F8").Value = "F<D<C<B<A<."
